$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 33 (row 22): mark Mon/Tue as done (copy the "completed" formatting
# used by the earlier weeks) and record 2 total days for that week.
$ws.Range("D9:E9").Copy()
$ws.Range("D22:E22").PasteSpecial(-4122)
$ws.Range("I22").Value = 2

# Move the active selection from A2 to B2.
$ws.Range("B2").Select()
